# Test case design cleanup: remove the "实际值" (actual value) and
# "执行结果" (execution result) test-output columns (F, G) for the data
# rows (2-20) — these were placeholder results, not part of the case
# design table. The header row (row 1) keeps its F1/G1 labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:G20").ClearContents()

# Move the active selection to where the last cleared cell used to be.
$ws.Range("F19").Select()
